# Weekly update: a new Cilantro price record (week of 2022-02-24) is added
# to the "Vega Monumental Concepción" sheet. All existing records from the
# old row 14 onward shift down by one pair of rows (Primera/Segunda) to make
# room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 14:15, pushing everything from the old
# row 14 down to row 16 (and so on, through the former last row 169,
# which ends up at row 171).
$ws.Rows("14:15").Insert()

# Row 14: "Primera" quality record for the new week.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = "2022-02-24"
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100112040
$ws.Range("G14").Value = "Cilantro"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 600
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = 650
$ws.Range("N14").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O14").Value = "Región de Ñuble"
$ws.Range("P14").Value = 650
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"

# Row 15: "Segunda" quality record for the same new week.
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = "2022-02-24"
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 100112040
$ws.Range("G15").Value = "Cilantro"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 500
$ws.Range("M15").Value = 500
$ws.Range("N15").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O15").Value = "Región de Ñuble"
$ws.Range("P15").Value = 500
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"
